# "sumo intro and conclusion"
#
# Turn on distinct odd/even headers & footers for the (only) section so
# Word splits the current single "default" header/footer into primary,
# even-page and first-page variants, then set their content:
#   - primary (odd/default) header: shortened title, no "with Dplyr"
#   - even-page header: emptied out (it inherits the old header1.xml part)
#   - first-page header/footer and the three footers: left blank
#
# WdHeaderFooterIndex constants: wdHeaderFooterPrimary = 1,
# wdHeaderFooterFirstPage = 2, wdHeaderFooterEvenPages = 3.
$wdHeaderFooterPrimary = 1
$wdHeaderFooterFirstPage = 2
$wdHeaderFooterEvenPages = 3

$d = $word.ActiveDocument
$section = $d.Sections(1)
$pageSetup = $section.PageSetup

# Enabling this splits off a distinct "even" header/footer from the
# existing "default" one, and (together with referencing the FirstPage
# header/footer below) causes Word to materialize header2/3.xml and
# footer1/2/3.xml with the correct rIds/ordering in sectPr.
$pageSetup.OddAndEvenPagesHeaderFooter = $true

$evenHeader = $section.Headers($wdHeaderFooterEvenPages)
$primaryHeader = $section.Headers($wdHeaderFooterPrimary)
$firstPageHeader = $section.Headers($wdHeaderFooterFirstPage)

$evenFooter = $section.Footers($wdHeaderFooterEvenPages)
$primaryFooter = $section.Footers($wdHeaderFooterPrimary)
$firstPageFooter = $section.Footers($wdHeaderFooterFirstPage)

# New shortened intro title in the primary (default/odd) header.
$primaryHeader.Range.Text = "Wrangling Sumo Wrestling Data"

# Everything else (even header, and all three footers) stays blank.
$evenFooter.Range.Text = ""
$primaryFooter.Range.Text = ""
$firstPageFooter.Range.Text = ""
$firstPageHeader.Range.Text = ""

# Turning this back off drops the now-unneeded <w:evenAndOddHeaders/>
# document setting while keeping the header/footer parts and the
# even/default/first headerReference & footerReference entries already
# written into sectPr.
$pageSetup.OddAndEvenPagesHeaderFooter = $false
